# acceptance: add AT-01 automation and close remaining partial rows
#
# WBS sheet: rows 50-62 ("Phase 3 - Traceability Engine" and its children)
# move from "Partial" to "Done", gain a Completed On date (column K) and
# five acceptance checkmarks (Schema / Validation / Permissions-Isolation /
# Workflow / Evidence -> columns L-P).
#
# Milestones sheet: M3 status flips from "In Progress" to "Done".

$wb = $excel.ActiveWorkbook
$wbs = $wb.Worksheets.Item("WBS")

$firstRow = 50
$lastRow = 62
$completedOn = "2026-04-17"
$check = "✅"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # H: Execution Status Partial -> Done
    $wbs.Cells.Item($r, 8).Value = "Done"

    # K: Completed On - text date, same display format as column J (s=51 / yyyy-mm-dd),
    # not an auto-converted Excel date serial. Force Text format before the write so
    # Excel doesn't coerce the date-shaped string into a real date serial, then restore
    # the yyyy-mm-dd display format (matches column J / I's styling).
    $wbs.Cells.Item($r, 11).NumberFormat = "@"
    $wbs.Cells.Item($r, 11).Value = $completedOn
    $wbs.Cells.Item($r, 11).NumberFormat = $wbs.Cells.Item($r, 10).NumberFormat

    # L-P: Schema / Validation / Permissions-Isolation / Workflow / Evidence checkmarks
    $wbs.Cells.Item($r, 12).Value = $check
    $wbs.Cells.Item($r, 13).Value = $check
    $wbs.Cells.Item($r, 14).Value = $check
    $wbs.Cells.Item($r, 15).Value = $check
    $wbs.Cells.Item($r, 16).Value = $check
}

$milestones = $wb.Worksheets.Item("Milestones")
$milestones.Range("F5").Value = "✅ Done"
